# Add 2022-Q1 sheet (new fund-holding snapshot) before the "总计" summary
# sheet, and add a corresponding row to "总计".

$wb = $excel.ActiveWorkbook

$refSheet   = $wb.Worksheets.Item("2021-Q4")
$totalSheet = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q1" worksheet right before "总计".
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($totalSheet)
$newSheet.Name = "2022-Q1"

# Borrow the header/row-label formatting used by the other quarterly
# sheets (bold, centered, bordered - style used for B1:H1 and column A).
$refSheet.Range("B1:H1").Copy()
$newSheet.Range("B1").PasteSpecial(-4122)
$refSheet.Range("A2").Copy()
$newSheet.Range("A2:A6").PasteSpecial(-4122)

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$newSheet.Range("A2").Value = 0
$newSheet.Range("A3").Value = 1
$newSheet.Range("A4").Value = 2
$newSheet.Range("A5").Value = 3
$newSheet.Range("A6").Value = 4

# Fund code / numeric-looking text columns must stay text (keep leading
# zeros / exact decimal strings), so force a text format before writing.
$newSheet.Range("B2:B6").NumberFormat = "@"
$newSheet.Range("D2:G6").NumberFormat = "@"

$newSheet.Range("B2").Value = "210003"
$newSheet.Range("C2").Value = "金鹰行业优势混合"
$newSheet.Range("D2").Value = "6.62"
$newSheet.Range("E2").Value = "88.02"
$newSheet.Range("F2").Value = "4.64"
$newSheet.Range("G2").Value = "0.3072"
$newSheet.Range("H2").Value = 4

$newSheet.Range("B3").Value = "013417"
$newSheet.Range("C3").Value = "博时核心资产精选混合A"
$newSheet.Range("D3").Value = "7.90"
$newSheet.Range("E3").Value = "78.99"
$newSheet.Range("F3").Value = "2.50"
$newSheet.Range("G3").Value = "0.1975"
$newSheet.Range("H3").Value = 9

$newSheet.Range("B4").Value = "005265"
$newSheet.Range("C4").Value = "博时厚泽回报灵活配置混合A"
$newSheet.Range("D4").Value = "2.19"
$newSheet.Range("E4").Value = "91.85"
$newSheet.Range("F4").Value = "3.00"
$newSheet.Range("G4").Value = "0.0657"
$newSheet.Range("H4").Value = 9

$newSheet.Range("B5").Value = "005266"
$newSheet.Range("C5").Value = "博时厚泽回报灵活配置混合C"
$newSheet.Range("D5").Value = "0.64"
$newSheet.Range("E5").Value = "91.85"
$newSheet.Range("F5").Value = "3.00"
$newSheet.Range("G5").Value = "0.0192"
$newSheet.Range("H5").Value = 9

$newSheet.Range("B6").Value = "013418"
$newSheet.Range("C6").Value = "博时核心资产精选混合C"
$newSheet.Range("D6").Value = "0.40"
$newSheet.Range("E6").Value = "78.99"
$newSheet.Range("F6").Value = "2.50"
$newSheet.Range("G6").Value = "0.0100"
$newSheet.Range("H6").Value = 9

# ---------------------------------------------------------------------
# 2. Add the 2022-Q1 roll-up row to "总计", above the existing rows.
#    NOTE: re-resolve "总计" by name here - the handle captured before
#    Worksheets.Add() above is a positional reference and, after the
#    insert shifted "总计" over by one slot, would otherwise point at
#    the freshly-added "2022-Q1" sheet instead.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 5
$totalSheet.Range("D2").Value = 0.6
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2

# Re-apply the original formatting that "Insert" does not cleanly carry
# over, by borrowing it from the (now renumbered, untouched) rows below.
$totalSheet.Range("A4").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
$totalSheet.Range("B4:D4").Copy()
$totalSheet.Range("B2:D2").PasteSpecial(-4122)

Write-Host "2022-Q1 sheet added and 总计 updated"
